$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension/measure metadata for the curated columns
$ws.Range("F2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("O2").Value = "iaest-measure:tipo-de-presupuesto"

$ws.Range("M3").Value = "dim"
$ws.Range("O3").Value = "medida"

$ws.Range("M4").Value = "URI-Municipio"

# Remove the now-obsolete mapping file references
$ws.Range("F5").Clear()
$ws.Range("O5").Clear()
